$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Status for "Making csv_operations module." moved from "In progress" to "Done"
$ws.Range("C3").Value = "Done"

# New row 6: PA3 task, Rex Liner, Done
$ws.Range("A6").Value = 45625
$ws.Range("B6").Value = "PA3"
$ws.Range("C6").Value = "Done"
$ws.Range("D6").Value = "Rex Liner"

# New row 7: final submission task, Drew Hutchinson, Done
$ws.Range("A7").Value = 45630
$ws.Range("B7").Value = "Fixing code and prep for final submission."
$ws.Range("C7").Value = "Done"
$ws.Range("D7").Value = "Drew Hutchinson"

# Match the date-cell formatting used by the rest of column A (avoid minting a new style)
$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C7").Select()
